$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 4794
$ws.Range("F10").Value = 1185
$ws.Range("F11").Value = 1627
$ws.Range("F12").Value = 852
$ws.Range("F13").Value = 528
$ws.Range("F14").Value = 2040
$ws.Range("F15").Value = 642
$ws.Range("F16").Value = 527
$ws.Range("F19").Value = 259
$ws.Range("F20").Value = 130
$ws.Range("F21").Value = 130
$ws.Range("F24").Value = 650
$ws.Range("F25").Value = 2555
$ws.Range("F29").Value = 1671
$ws.Range("F34").Value = 4415

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 401
$ws.Range("F6").Value = 4180
$ws.Range("F16").Value = 16
$ws.Range("F25").Value = 14
$ws.Range("F26").Value = 211

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1388
$ws.Range("F5").Value = 1760
$ws.Range("F7").Value = 435

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1388
$ws.Range("F4").Value = 1760
$ws.Range("F6").Value = 435
$ws.Range("F8").Value = 4794
$ws.Range("F15").Value = 1185
$ws.Range("F16").Value = 1627
$ws.Range("F20").Value = 852
$ws.Range("F21").Value = 528
$ws.Range("F22").Value = 2040
$ws.Range("F23").Value = 642
$ws.Range("F24").Value = 527
$ws.Range("F27").Value = 259
$ws.Range("F29").Value = 130
$ws.Range("F30").Value = 130
$ws.Range("F35").Value = 650
$ws.Range("F38").Value = 2555
$ws.Range("F44").Value = 1671
$ws.Range("F49").Value = 4415
